$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 22:20"

# --- Update row 4 (Estados Unidos) ---
$ws.Range("B4").Value = 270062
$ws.Range("C4").Value = 25185
$ws.Range("E4").Value = 251120
$ws.Range("G4").Value = 857
$ws.Range("H4").Value = 6927

# --- Update row 6 (España) ---
$ws.Range("E6").Value = 76188
$ws.Range("G6").Value = 661
$ws.Range("H6").Value = 11009

# --- Update row 7 (Alemania) ---
$ws.Range("B7").Value = 91159
$ws.Range("C7").Value = 6365
$ws.Range("E7").Value = 65309
$ws.Range("G7").Value = 168
$ws.Range("H7").Value = 1275

# --- Update row 44 (Peru) ---
$ws.Range("E44").Value = 997
$ws.Range("G44").Value = 6
$ws.Range("H44").Value = 61

# --- Update row 48 (Serbia) ---
$ws.Range("D48").Value = 54
$ws.Range("E48").Value = 1383

# --- Costa de Marfil moves up in the ranking (new case counts push it
#     above Senegal/Ghana/Malta). Shift those three rows down one slot
#     and give Costa de Marfil its updated totals at row 100. ---
$ws.Range("A103").Value = "Malta"
$ws.Range("B103").Value = 202
$ws.Range("C103").Value = 6
$ws.Range("D103").Value = 2
$ws.Range("E103").Value = 200
$ws.Range("F103").Value = 2
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0

$ws.Range("A102").Value = "Ghana"
$ws.Range("B102").Value = 205
$ws.Range("C102").Value = 1
$ws.Range("D102").Value = 31
$ws.Range("E102").Value = 169
$ws.Range("F102").Value = 2
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 5

$ws.Range("A101").Value = "Senegal"
$ws.Range("B101").Value = 207
$ws.Range("C101").Value = 12
$ws.Range("D101").Value = 66
$ws.Range("E101").Value = 140
$ws.Range("F101").Value = 1
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 1

$ws.Range("A100").Value = "Costa de Marfil"
$ws.Range("B100").Value = 218
$ws.Range("C100").Value = 24
$ws.Range("D100").Value = 19
$ws.Range("E100").Value = 198
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 1
